# Add handling of transcript files to BatchIngest
#
# The manifest sheet gains a new 3-column "Transcript" field group
# (Transcript File / Transcript Label / Machine Generated), inserted right
# after the existing "Caption Language" column (which ends at column V),
# pushing every later column group (the plain File/Offset/Label group, the
# Skip Transcoding/Absolute Location/Date Digitized group, the second
# File/Label/Skip Transcoding group and the Abstract/Statement Of
# Responsibility columns) three places to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 3 blank columns before the old column W ("File" of the
#    "Absolute" group). Excel shifts all data/formatting/column-widths that
#    lived at W onward three columns to the right (-> Z onward).
$ws.Range("W1:Y1").EntireColumn.Insert()

# 2. New header row (row 2) labels for the freshly inserted columns.
$ws.Range("W2").Value = "Transcript File"
$ws.Range("X2").Value = "Transcript Label"
$ws.Range("Y2").Value = "Machine Generated"

# 3. New data row (row 3) values for the transcript of the Sheephead
#    Mountain asset - mirrors the existing caption row right next to it.
$ws.Range("W3").Value = $ws.Range("T3").Value()
$ws.Range("X3").Value = "Sheephead Transcript"
$ws.Range("Y3").Value = "yes"

# Rows 4 and 5 have nothing in the new Transcript columns, matching the
# blanks that come from the column insert.

# 4. The hyperlink that decorated the old AA3 ("Absolute Location" of the
#    first file group) now lives three columns over, at AD3. The insert
#    shifts the cell's text but not the worksheet's Hyperlinks collection,
#    so re-anchor it explicitly.
$oldHyperlinkCell = $ws.Range("AA3")
$hyperlinkAddress = "file:///tmp/sheephead_mountain_master.mov"
$newHyperlinkCell = $ws.Range("AD3")
$oldHyperlinkCell.Hyperlinks.Delete()
$ws.Hyperlinks.Add($newHyperlinkCell, $hyperlinkAddress, "", "", $newHyperlinkCell.Value())

# 5. Reflect the edit's final cursor position/selection.
$ws.Range("X3").Select()
